$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cell A1
$ws.Range("A1").Value = "Cluster Name"

# Update the "3824 Estia Health South Morang" Active cases value (row 8) from 48 to 52
$ws.Cells.Item(8, 2).Value = 52

# Delete rows for clusters that are no longer reported (in descending order
# so earlier row numbers remain valid while deleting)
$rowsToDelete = @(46, 45, 43, 40, 35, 33, 12)
foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}
